$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 5494937
$ws.Range("I58").Value = 6493689.5
$ws.Range("J58").Value = 1800
$ws.Range("K58").Value = 19481068.5
$ws.Range("L58").Value = 5400
$ws.Range("M58").Value = -19480918.5
$ws.Range("N58").Value = -5700
$ws.Range("H64").Value = 3062.0688
$ws.Range("I64").Value = 2926.3157
$ws.Range("J64").Value = 3320
$ws.Range("K64").Value = 2926.3157
$ws.Range("L64").Value = 3320
$ws.Range("M64").Value = -2678.3157
$ws.Range("N64").Value = -3816
$ws.Range("H67").Value = 3062.0688
$ws.Range("I67").Value = 2926.3157
$ws.Range("J67").Value = 3320
$ws.Range("K67").Value = 2926.3157
$ws.Range("L67").Value = 3320
$ws.Range("M67").Value = -2068.3157
$ws.Range("N67").Value = -5036
$ws.Range("H76").Value = 3003
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 3003
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H98").Value = 36196.637
$ws.Range("I98").Value = 883.2857
$ws.Range("J98").Value = 777777
$ws.Range("K98").Value = 883.2857
$ws.Range("L98").Value = 777777
$ws.Range("M98").Value = 614.7143
$ws.Range("N98").Value = -780773
$ws.Range("H107").Value = 4429.385
$ws.Range("I107").Value = 4798.522
$ws.Range("J107").Value = 1599.3334
$ws.Range("K107").Value = 4798.522
$ws.Range("L107").Value = 1599.3334
$ws.Range("M107").Value = -2878.522
$ws.Range("N107").Value = -5439.3334
$ws.Range("H122").Value = 36196.637
$ws.Range("I122").Value = 883.2857
$ws.Range("J122").Value = 777777
$ws.Range("K122").Value = 2649.8571
$ws.Range("L122").Value = 2333331
$ws.Range("M122").Value = -199.8571000000002
$ws.Range("N122").Value = -2338231
$ws.Range("H125").Value = 1621.4
$ws.Range("I125").Value = 1463.5
$ws.Range("J125").Value = 1726.6666
$ws.Range("K125").Value = 13171.5
$ws.Range("L125").Value = 15539.9994
$ws.Range("M125").Value = -10711.5
$ws.Range("N125").Value = -20459.9994
$ws.Range("H127").Value = 994.6
$ws.Range("I127").Value = 256.6
$ws.Range("J127").Value = 1732.6
$ws.Range("K127").Value = 769.8000000000001
$ws.Range("L127").Value = 5197.799999999999
$ws.Range("M127").Value = 4190.2
$ws.Range("N127").Value = -15117.8
$ws.Range("H129").Value = 903.37036
$ws.Range("J129").Value = 904.1667
$ws.Range("L129").Value = 2712.5001
$ws.Range("N129").Value = -12712.5001

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7687.073
$ws.Range("I32").Value = 5600.6562
$ws.Range("J32").Value = 15105.444
$ws.Range("K32").Value = 5600.6562
$ws.Range("L32").Value = 15105.444
$ws.Range("M32").Value = -5313.6562
$ws.Range("N32").Value = -15679.444
$ws.Range("H63").Value = 3421.5625
$ws.Range("I63").Value = 2155.5
$ws.Range("J63").Value = 5531.6665
$ws.Range("K63").Value = 2155.5
$ws.Range("L63").Value = 5531.6665
$ws.Range("M63").Value = -1469.5
$ws.Range("N63").Value = -6903.6665
$ws.Range("H66").Value = 3421.5625
$ws.Range("I66").Value = 2155.5
$ws.Range("J66").Value = 5531.6665
$ws.Range("K66").Value = 10777.5
$ws.Range("L66").Value = 27658.3325
$ws.Range("M66").Value = -7345.5
$ws.Range("N66").Value = -34522.3325

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7097536
$ws.Range("I31").Value = 1824.421
$ws.Range("J31").Value = 11912483
$ws.Range("K31").Value = 1824.421
$ws.Range("L31").Value = 11912483
$ws.Range("M31").Value = -1529.421
$ws.Range("N31").Value = -11913073
$ws.Range("H34").Value = 7097536
$ws.Range("I34").Value = 1824.421
$ws.Range("J34").Value = 11912483
$ws.Range("K34").Value = 1824.421
$ws.Range("L34").Value = 11912483
$ws.Range("M34").Value = -1622.421
$ws.Range("N34").Value = -11912887
$ws.Range("H134").Value = 535940.5600000001
$ws.Range("I134").Value = 641662.1
$ws.Range("J134").Value = 203672.86
$ws.Range("K134").Value = 1924986.3
$ws.Range("L134").Value = 611018.58
$ws.Range("M134").Value = -1922451.3
$ws.Range("N134").Value = -616088.58

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 330.51614
$ws.Range("J12").Value = 489.5
$ws.Range("L12").Value = 1468.5
$ws.Range("N12").Value = -1814.5
$ws.Range("H20").Value = 371.42856
$ws.Range("I20").Value = 371.42856
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1114.28568
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -887.28568
$ws.Range("N20").ClearContents()
$ws.Range("H32").Value = 730
$ws.Range("J32").Value = 730
$ws.Range("L32").Value = 2190
$ws.Range("N32").Value = -2756
$ws.Range("H33").Value = 14050085
$ws.Range("I33").Value = 66
$ws.Range("J33").Value = 19067950
$ws.Range("K33").Value = 396
$ws.Range("L33").Value = 114407700
$ws.Range("M33").Value = -113
$ws.Range("N33").Value = -114408266
$ws.Range("H34").Value = 2922.1667
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 2922.1667
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 8766.500100000001
$ws.Range("N34").Value = -8934.500100000001
$ws.Range("M34").ClearContents()
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").ClearContents()
$ws.Range("H129").Value = 178298.83
$ws.Range("I129").Value = 1001010
$ws.Range("J129").Value = 2003.5714
$ws.Range("K129").Value = 3003030
$ws.Range("L129").Value = 6010.7142
$ws.Range("M129").Value = -2998030
$ws.Range("N129").Value = -16010.7142
$ws.Range("H139").Value = 99925.25999999999
$ws.Range("I139").Value = 378583
$ws.Range("J139").Value = 3000.8262
$ws.Range("K139").Value = 1135749
$ws.Range("L139").Value = 9002.4786
$ws.Range("M139").Value = -1130609
$ws.Range("N139").Value = -19282.4786

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1219
$ws.Range("I102").Value = 1500
$ws.Range("J102").Value = 1078.5
$ws.Range("K102").Value = 1500
$ws.Range("L102").Value = 1078.5
$ws.Range("M102").Value = 122
$ws.Range("N102").Value = -4322.5

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2744
$ws.Range("I40").Value = 2804.125
$ws.Range("J40").Value = 2503.5
$ws.Range("K40").Value = 2804.125
$ws.Range("L40").Value = 2503.5
$ws.Range("M40").Value = -2668.125
$ws.Range("N40").Value = -2775.5
$ws.Range("H61").Value = 4062.3076
$ws.Range("I61").Value = 4061
$ws.Range("J61").Value = 4066.6667
$ws.Range("K61").Value = 4061
$ws.Range("L61").Value = 4066.6667
$ws.Range("M61").Value = -3859
$ws.Range("N61").Value = -4470.6667
$ws.Range("H113").Value = 4062.3076
$ws.Range("I113").Value = 4061
$ws.Range("J113").Value = 4066.6667
$ws.Range("K113").Value = 4061
$ws.Range("L113").Value = 4066.6667
$ws.Range("M113").Value = -1891
$ws.Range("N113").Value = -8406.6667

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4480.4
$ws.Range("I62").Value = 6701
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 6701
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -6077
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 4480.4
$ws.Range("I65").Value = 6701
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 33505
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -30385
$ws.Range("N65").Value = -21240
$ws.Range("H122").Value = 1682310.6
$ws.Range("I122").Value = 1906418.8
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 5719256.4
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -5716806.4
$ws.Range("N122").Value = -9400
